$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4683
$ws1.Range("F3").Value = 2723
$ws1.Range("F5").Value = 2741
$ws1.Range("F9").Value = 1719
$ws1.Range("F10").Value = 743
$ws1.Range("F12").Value = 212
$ws1.Range("F15").Value = 296
$ws1.Range("G18").Value = 54
$ws1.Range("G19").Value = 54
$ws1.Range("F22").Value = 644
$ws1.Range("F23").Value = 736
$ws1.Range("F28").Value = 1460
$ws1.Range("F29").Value = 316
$ws1.Range("F31").Value = 1423
$ws1.Range("F32").Value = 2285
$ws1.Range("F33").Value = 381
$ws1.Range("F37").Value = 55
$ws1.Range("F39").Value = 772
$ws1.Range("F40").Value = 1456
$ws1.Range("F41").Value = 197
$ws1.Range("F44").Value = 22

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 33

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4683
$ws4.Range("F3").Value = 2723
$ws4.Range("F4").Value = 2741
$ws4.Range("F5").Value = 1719
$ws4.Range("F8").Value = 743
$ws4.Range("F10").Value = 212
$ws4.Range("F13").Value = 296
$ws4.Range("G16").Value = 54
$ws4.Range("G17").Value = 54
$ws4.Range("F19").Value = 644
$ws4.Range("F20").Value = 736
$ws4.Range("F28").Value = 1460
$ws4.Range("F29").Value = 316
$ws4.Range("F33").Value = 2285
$ws4.Range("F34").Value = 381
$ws4.Range("F38").Value = 33
$ws4.Range("F41").Value = 55
$ws4.Range("F43").Value = 772
$ws4.Range("F44").Value = 1456
$ws4.Range("F46").Value = 197

